# Updates the crypto price/volume snapshot on Sheet1 (columns D = Price,
# E = Volume(1h)) to the latest scraped values, mirroring the GitHub
# Actions "Updated cryptos list" commit.
#
# Some new Price values are plain numeric-looking strings (e.g. "7.63",
# "0.515", "460.80"). Excel's Range.Value setter auto-coerces such text to
# a real number (losing the original text formatting / trailing zeros and
# changing the stored cell type away from a string). To keep those cells
# as literal text - matching the workbook's existing inline-string cells -
# we briefly force the cell to Text number-format before assigning the
# value, then restore the original "Normal" style so no visible formatting
# changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    $cell = $ws.Range($rangeAddress)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# row 2 - Bitcoin
Set-TextValue "D2" '62.909.33'
$ws.Range("E2").Value = '  -1.44%  '

# row 3 - Ethereum
Set-TextValue "D3" '3.164.53'
$ws.Range("E3").Value = '  -4.75%  '

# row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.02%  '

# row 5 - BNB
Set-TextValue "D5" '589.67'
$ws.Range("E5").Value = '  -2.41%  '

# row 6 - Solana
Set-TextValue "D6" '133.91'
$ws.Range("E6").Value = '  -6.18%  '

# row 8 - LidoStakedEther
Set-TextValue "D8" '3.162.45'
$ws.Range("E8").Value = '  -4.79%  '

# row 9 - XRP
Set-TextValue "D9" '0.515'
$ws.Range("E9").Value = '  -0.88%  '

# row 10 - Dogecoin
$ws.Range("E10").Value = '  -6.41%  '

# row 11 - Toncoin
$ws.Range("E11").Value = '  -5.53%  '

# row 12 - Cardano
$ws.Range("E12").Value = '  -3.43%  '

# row 13 - ShibaInu
Set-TextValue "D13" '0.0000235'
$ws.Range("E13").Value = '  -5.13%  '

# row 14 - Avalanche
Set-TextValue "D14" '34.87'
$ws.Range("E14").Value = '  -0.65%  '

# row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" '3.685.25'
$ws.Range("E15").Value = '  -4.82%  '

# row 16 - TRON
$ws.Range("E16").Value = '  -1.14%  '

# row 17 - WrappedEther
Set-TextValue "D17" '3.168.94'
$ws.Range("E17").Value = '  -4.73%  '

# row 18 - WrappedBTC
Set-TextValue "D18" '62.896.29'
$ws.Range("E18").Value = '  -1.57%  '

# row 19 - Polkadot
$ws.Range("E19").Value = '  -4.55%  '

# row 20 - BitcoinCash
Set-TextValue "D20" '460.80'
$ws.Range("E20").Value = '  -4.00%  '

# row 21 - Chainlink
Set-TextValue "D21" '13.87'
$ws.Range("E21").Value = '  -1.84%  '

# row 22 - Polygon
$ws.Range("E22").Value = '  -5.72%  '

# row 23 - Uniswap
Set-TextValue "D23" '7.63'
$ws.Range("E23").Value = '  -4.76%  '

# row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" '13.39'
$ws.Range("E24").Value = '  -3.44%  '

# row 25 - Litecoin (Price changes; Volume(1h) unchanged)
Set-TextValue "D25" '82.96'

# row 26 - Dai
$ws.Range("E26").Value = '  -0.08%  '

# row 27 - FirstDigitalUSD
$ws.Range("E27").Value = '  +0.04%  '

# row 29 - NEARProtocol
Set-TextValue "D29" '6.80'
$ws.Range("E29").Value = '  -5.35%  '

# row 30 - RenderToken
Set-TextValue "D30" '7.71'
$ws.Range("E30").Value = '  -6.62%  '

# row 31 - ImmutableX
$ws.Range("E31").Value = '  -6.00%  '

# row 32 - EthereumClassic (Price changes; Volume(1h) unchanged)
Set-TextValue "D32" '27.08'

# row 33 - Hedera
$ws.Range("E33").Value = '  -4.59%  '

# row 34 - Stacks
Set-TextValue "D34" '2.35'
$ws.Range("E34").Value = '  -7.07%  '

# row 35 - Mantle
Set-TextValue "D35" '1.03'
$ws.Range("E35").Value = '  -6.23%  '

# row 36 - Filecoin
Set-TextValue "D36" '5.80'
$ws.Range("E36").Value = '  -4.50%  '

# row 37 - OKB
Set-TextValue "D37" '51.28'
$ws.Range("E37").Value = '  -2.18%  '

# row 38 - PEPE
$ws.Range("E38").Value = '  -6.66%  '

# row 39 - VeChain
$ws.Range("E39").Value = '  -3.21%  '

# row 40 - Bittensor
Set-TextValue "D40" '401.47'
$ws.Range("E40").Value = '  -6.83%  '

# row 41 - Cosmos
Set-TextValue "D41" '8.08'
$ws.Range("E41").Value = '  -3.15%  '

# row 42 - Kaspa
$ws.Range("E42").Value = '  -4.23%  '

# row 43 - dogwifhat
$ws.Range("E43").Value = '  -5.61%  '

# row 44 - Maker
$ws.Range("D44").Value = '2.790.30'
$ws.Range("E44").Value = '  -10.84%  '

# row 45 - TheGraph
$ws.Range("E45").Value = '  -6.42%  '

# row 46 - USDe
Set-TextValue "D46" '0.999'
$ws.Range("E46").Value = '  +0.04%  '

# row 47 - Fetch.AI
$ws.Range("E47").Value = '  -6.04%  '

# row 48 - Monero
Set-TextValue "D48" '124.79'
$ws.Range("E48").Value = '  +0.50%  '

# row 49 - InjectiveProtocol
Set-TextValue "D49" '25.13'
$ws.Range("E49").Value = '  -4.91%  '

# row 50 - Arweave
Set-TextValue "D50" '34.25'
$ws.Range("E50").Value = '  -5.34%  '

# row 51 - Stellar
$ws.Range("E51").Value = '  -2.35%  '
